$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(11, 10, 10, 10, 9, 9, 6, 6, 6, 6, 6, 6, 5, 2)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B17").Select()
